$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$docente = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# The workbook originally had a standalone row (no column-A label) right after
# "Docentes responsáveis:" holding just the professor's name in B/C. That row
# is removed, shifting every row below it up by one.
$ws.Rows(13).Delete()

# --- Content edits (row numbers below are POST-deletion positions) ---

# "Objetivos:" row: the long objectives paragraph is replaced by the docente's
# identification (reusing the text that used to sit in the deleted row).
$ws.Range("B10").Value = $docente
$ws.Range("C10").Value = $docente

# "Programa resumido:" row (old row 14): long description replaced with "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" row (old row 16): long syllabus text replaced with the activation
# date text "01/01/2020" (kept as literal text, not an auto-converted date).
$ws.Range("B15").Value = "'01/01/2020"
$ws.Range("C15").Value = "'01/01/2020"

# "Método:" row (old row 19): value replaced with the docente's identification.
$ws.Range("B18").Value = $docente
$ws.Range("C18").Value = $docente

# "Critério:" row (old row 20): now holds the text that used to be the
# "Método:" evaluation description.
$ws.Range("B19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."
$ws.Range("C19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."

# "Norma de recuperação:" row (old row 21): now holds the text that used to be
# the "Critério:" description.
$ws.Range("B20").Value = "Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois."
$ws.Range("C20").Value = "Serão aplicadas duas avaliações para compor a média que será a soma das duas provas, sendo o resultado dividido por dois."

# "Bibliografia:" row (old row 22): now holds the text that used to be the
# "Norma de recuperação:" description (the huge bibliography text is gone).
$ws.Range("B21").Value = "A Nota Final será composta pela Média obtida da Nota do Período somada à Nota de Recuperação e dividido por dois"
$ws.Range("C21").Value = "A Nota Final será composta pela Média obtida da Nota do Período somada à Nota de Recuperação e dividido por dois"
